# Updates "Projektkosten PAMS (Zentral)" worksheet:
#  - Expands the task-description texts for the "interne Leistungen" rows
#    (16-20) with more detailed/longer descriptions.
#  - Applies the formatting that goes along with the longer, wrapped text
#    (top-vertical alignment on the row, wrapped description column, taller
#    rows for the now multi-line descriptions).
#  - Moves the active selection as it was left after editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the task descriptions (column C) for rows 16-20 ---------
$ws.Range("C16").Value = "25h, Logik, Ausgabe, Eingabe, Datenbankzugriff, Login(Verschlüsselung), "
$ws.Range("C17").Value = "8h, DB-Design, Erstellung"
$ws.Range("C18").Value = "40h, Konzeption, Dokumenation, Aufgabenverteilung, Budgetverwaltung, Meilensteine festlegen"
$ws.Range("C20").Value = "32h, Machbarkeit, Umfeld, Risiko, Statusberichte, Kommunikationsrichtlinien, Dokumentationsrichtlinien"

# --- 2. Formatting: column A (occupation) gets top-aligned text ---------
$ws.Range("A16:A20").VerticalAlignment = -4160   # xlTop

# --- 3. Formatting: column B (cost) gets top-aligned numbers ------------
$ws.Range("B16:B20").VerticalAlignment = -4160   # xlTop

# --- 4. Formatting: column C wraps for the rows with longer text -------
$ws.Range("C16").WrapText = $true
$ws.Range("C18").WrapText = $true
$ws.Range("C20").WrapText = $true

# --- 5. Row heights grow to fit the now-wrapped, longer descriptions ---
$ws.Rows(16).RowHeight = 30
$ws.Rows(18).RowHeight = 30
$ws.Rows(20).RowHeight = 30

# --- 6. Leave the selection where the author left it after editing -----
$ws.Range("G15").Select()
